# Automatische test-sync: 2025-06-20 10:00:50
# Adds the new incoming mail log entry to the "Logs" sheet and refreshes
# the category totals on the "Dashboard" sheet accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append the new mail log row (row 6)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Vragen over samenwerking"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D6").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F6").Value = "2025-06-20 10:00:12"
$logs.Range("G6").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too.
$logs.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D6"))
$logs.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G6"))

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: re-sort the category summary table, the new row
#    raised "Samenwerking / Partnerverzoek" to 2 occurrences.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Samenwerking / Partnerverzoek"
$dashboard.Range("B2").Value = 2
$dashboard.Range("A3").Value = "Afmelding / Nieuwsbrief"
$dashboard.Range("B3").Value = 1
$dashboard.Range("A4").Value = "Productinformatie"
$dashboard.Range("B4").Value = 1
